# WCP.docx edit: replace the inline "spacer" picture (the one linked to
# rId23 / media/rId23.png) that sits in the "Submission of WCP" section
# with a hyperlink whose visible text is the image's original URL.
#
#   <w:r><w:drawing>...pic r:embed="rId23".../></w:drawing></w:r>
# becomes
#   <w:hyperlink r:id="rId23">
#     <w:r><w:rPr><w:rStyle w:val="Hyperlink"/></w:rPr>
#       <w:t>https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/Others/WCP.jpg</w:t>
#     </w:r>
#   </w:hyperlink>

$d = $word.ActiveDocument

$url = "https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/Others/WCP.jpg"

# Locate the lone inline picture in the document (the 1x1px image that was
# embedded via rId23) and grab its Range before removing it.
$shp = $d.InlineShapes.Item(1)
$rng = $shp.Range
$rng.Select()
$shp.Delete()

# Turn that now-empty range into a hyperlink, using the URL both as the
# target address and as the displayed text.
$d.Hyperlinks.Add($rng, $url, "", "", $url) | Out-Null
